# Update the MPA test automation upload file values on the "Data" sheet.
# K column: 60000178 -> 60000215
# L column: 133      -> 165
# N column: 60000179 -> 60000216
# O column: 134      -> 166

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$kRows = @(6, 9, 11, 14, 16, 19, 21, 24, 26, 29)
foreach ($r in $kRows) {
    $ws.Range("K$r").Value = 60000215
}

$lRows = @(7, 8, 10, 12, 13, 15, 17, 18, 20, 22, 23, 25, 27, 28)
foreach ($r in $lRows) {
    $ws.Range("L$r").Value = 165
}

$nRows = @(7, 11, 12, 16, 17, 21, 22, 26, 27)
foreach ($r in $nRows) {
    $ws.Range("N$r").Value = 60000216
}

$oRows = @(8, 13, 18, 23, 28)
foreach ($r in $oRows) {
    $ws.Range("O$r").Value = 166
}
